$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Remove the gray highlight fill from A66:A69 (array_01..array_04), which
#    were previously marked as the "latest batch" of added rows. A new batch
#    (rows 127/128) is being added now, so these lose their highlight.
# ---------------------------------------------------------------------------
$ws.Range("A66:A69").Interior.Pattern = -4142   # xlPatternNone ("No Fill")

# ---------------------------------------------------------------------------
# 2) Add row 127: insert_028 (boolean field insert decimal test case)
# ---------------------------------------------------------------------------
$ws.Range("A126:K126").Copy()
$ws.Range("A127:K127").PasteSpecial(-4122)      # xlPasteFormats
$ws.Range("O126").Copy()
$ws.Range("O127").PasteSpecial(-4122)           # xlPasteFormats

$ws.Range("A127").Value = "insert_028"
$ws.Range("B127").Value = "y"
$ws.Range("C127").Value = "布尔型字段插入小数"
$ws.Range("D127").Value = "insert"
$ws.Range("F127").Value = "schema1"
$ws.Range("G127").Value = "insert_value24"
$ws.Range("H127").Value = "4"
$ws.Range("I127").Value = "select * from `$schema1"
$ws.Range("J127").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/dml/insert/expectedresult/insert_028.csv"
$ws.Range("O127").Value = "csv_containsAll"

# ---------------------------------------------------------------------------
# 3) Add row 128: array_05 (boolean array insert decimal element test case)
# ---------------------------------------------------------------------------
$ws.Range("A126:K126").Copy()
$ws.Range("A128:K128").PasteSpecial(-4122)      # xlPasteFormats
$ws.Range("O126").Copy()
$ws.Range("O128").PasteSpecial(-4122)           # xlPasteFormats

$ws.Range("A128").Value = "array_05"
$ws.Range("B128").Value = "y"
$ws.Range("C128").Value = "布尔型数组插入元素为小数"
$ws.Range("D128").Value = "ComplexDataType"
$ws.Range("E128").Value = "Array"
$ws.Range("F128").Value = "array10"
$ws.Range("G128").Value = "array10_value58"
$ws.Range("H128").Value = "3"
$ws.Range("I128").Value = "select in_use from `$array10"
$ws.Range("J128").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/dml/insert/expectedresult/complexdatatype/array/array_005.csv"
$ws.Range("O128").Value = "csv_containsAll"

# A128 (the new TestID cell) carries the same "no fill" style that A66:A69
# now use (copy the already-cleared format from A66, which preserves the
# distinct style/fill table entry rather than re-deriving a deduped one).
$ws.Range("A66").Copy()
$ws.Range("A128").PasteSpecial(-4122)           # xlPasteFormats

# ---------------------------------------------------------------------------
# 4) Update the view: selection moved to B118, scrolled near the new rows.
# ---------------------------------------------------------------------------
$excel.Windows.Item(1).ScrollRow = 73
$ws.Range("B118").Select()
